$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry per-row data which shifts down by one row: D,J,K,L,M,P
# (= 4,10,11,12,13,16). Every other column (A,B,C,E,F,G,H,I,N,O,Q,R) is
# constant across the whole data block, so the newly-inserted row can just
# copy row 138's values for those.
$shiftCols = @(4, 10, 11, 12, 13, 16)
$lastRow = 138
$newRow = $lastRow + 1

# Shift rows 42..138 down into rows 43..139 (process bottom-up so a source
# row is never clobbered before it has been read). Rows 43..138 already
# exist with the correct per-column styles, so only the value needs to move.
# Row 139 is brand new, so its style must be copied explicitly (it starts
# out as a plain/default-styled cell).
for ($r = $lastRow; $r -ge 42; $r--) {
    $dst = $r + 1
    foreach ($c in $shiftCols) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($dst, $c)
        $dstCell.Value2 = $srcCell.Value2
        # Only column D (4) carries an explicit date style; the new row's
        # cell starts out default-styled, so copy it across there. The
        # other shifted columns (J,K,L,M,P) use the default style already,
        # both before and after, so leave them alone.
        if ($dst -eq $newRow -and $c -eq 4) {
            $dstCell.NumberFormat = $srcCell.NumberFormat
        }
    }
}

# Row 139 is brand new - copy the constant columns from row 138 too. These
# cells use the sheet's default (General) style already, so no style work
# is needed here.
$staticCols = @(1, 2, 3, 5, 6, 7, 8, 9, 14, 15, 17, 18)
foreach ($c in $staticCols) {
    $srcCell = $ws.Cells.Item($lastRow, $c)
    $dstCell = $ws.Cells.Item($newRow, $c)
    $dstCell.Value2 = $srcCell.Value2
}

# Row 42 now gets the brand-new observation (new date, new volume); the
# min/max/avg price columns (K,L,M,P) are unchanged from the original row 42.
$ws.Cells.Item(42, 4).Value2 = 44536
$ws.Cells.Item(42, 10).Value2 = 2700
